$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("J10").Value = 1.05
$ws.Range("K10").Value = 11

# Row 12
$ws.Range("G12").Value = 1.8
$ws.Range("H12").Value = 3.3
$ws.Range("I12").Value = 3.8
$ws.Range("J12").Value = 1.05
$ws.Range("K12").Value = 11
$ws.Range("L12").Value = 1.25
$ws.Range("M12").Value = 3.75
$ws.Range("N12").Value = 1.85
$ws.Range("O12").Value = 1.95
$ws.Range("P12").Value = 1.36
$ws.Range("Q12").Value = 3
$ws.Range("R12").Value = 1.73
$ws.Range("S12").Value = 2
$ws.Range("T12").Value = 8
$ws.Range("U12").Value = 9.5
$ws.Range("V12").Value = 8.5
$ws.Range("W12").Value = 15
$ws.Range("X12").Value = 15
$ws.Range("Y12").Value = 26
$ws.Range("Z12").Value = 11
$ws.Range("AA12").Value = 7
$ws.Range("AB12").Value = 15
$ws.Range("AC12").Value = 41
$ws.Range("AD12").Value = 201
$ws.Range("AE12").Value = 12
$ws.Range("AF12").Value = 21
$ws.Range("AG12").Value = 13
$ws.Range("AH12").Value = 41
$ws.Range("AI12").Value = 29
$ws.Range("AJ12").Value = 34

# Row 13
$ws.Range("G13").Value = 2.05
$ws.Range("H13").Value = 3.25
$ws.Range("I13").Value = 3.1
$ws.Range("J13").Value = 1.05
$ws.Range("K13").Value = 11
$ws.Range("L13").Value = 1.29
$ws.Range("M13").Value = 3.5
$ws.Range("N13").Value = 1.9
$ws.Range("O13").Value = 1.9
$ws.Range("P13").Value = 1.36
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = 1.73
$ws.Range("S13").Value = 2
$ws.Range("T13").Value = 8
$ws.Range("U13").Value = 11
$ws.Range("V13").Value = 9
$ws.Range("W13").Value = 19
$ws.Range("X13").Value = 17
$ws.Range("Y13").Value = 26
$ws.Range("Z13").Value = 11
$ws.Range("AA13").Value = 6.5
$ws.Range("AB13").Value = 13
$ws.Range("AC13").Value = 41
$ws.Range("AD13").Value = 201
$ws.Range("AE13").Value = 11
$ws.Range("AF13").Value = 17
$ws.Range("AG13").Value = 12
$ws.Range("AH13").Value = 34
$ws.Range("AI13").Value = 26
$ws.Range("AJ13").Value = 34

# Row 14
$ws.Range("G14").Value = 1.55
$ws.Range("H14").Value = 3.6
$ws.Range("I14").Value = 5.25
$ws.Range("J14").Value = 1.06
$ws.Range("K14").Value = 10
$ws.Range("L14").Value = 1.29
$ws.Range("M14").Value = 3.5
$ws.Range("N14").Value = 1.85
$ws.Range("O14").Value = 1.95
$ws.Range("P14").Value = 1.36
$ws.Range("Q14").Value = 3
$ws.Range("R14").Value = 1.91
$ws.Range("S14").Value = 1.8
$ws.Range("T14").Value = 7
$ws.Range("U14").Value = 7
$ws.Range("V14").Value = 8.5
$ws.Range("W14").Value = 11
$ws.Range("X14").Value = 13
$ws.Range("Y14").Value = 29
$ws.Range("Z14").Value = 10
$ws.Range("AA14").Value = 7.5
$ws.Range("AB14").Value = 17
$ws.Range("AC14").Value = 51
$ws.Range("AD14").Value = 301
$ws.Range("AE14").Value = 15
$ws.Range("AF14").Value = 29
$ws.Range("AG14").Value = 19
$ws.Range("AH14").Value = 67
$ws.Range("AI14").Value = 41
$ws.Range("AJ14").Value = 41

# Row 15
$ws.Range("G15").Value = 1.8
$ws.Range("H15").Value = 3.4
$ws.Range("I15").Value = 3.8
$ws.Range("J15").Value = 1.05
$ws.Range("K15").Value = 11
$ws.Range("L15").Value = 1.29
$ws.Range("M15").Value = 3.5
$ws.Range("N15").Value = 1.9
$ws.Range("O15").Value = 1.9
$ws.Range("P15").Value = 1.36
$ws.Range("Q15").Value = 3
$ws.Range("R15").Value = 1.8
$ws.Range("S15").Value = 1.91
$ws.Range("T15").Value = 7.5
$ws.Range("U15").Value = 9
$ws.Range("V15").Value = 8.5
$ws.Range("W15").Value = 15
$ws.Range("X15").Value = 15
$ws.Range("Y15").Value = 26
$ws.Range("Z15").Value = 11
$ws.Range("AA15").Value = 7
$ws.Range("AB15").Value = 15
$ws.Range("AC15").Value = 51
$ws.Range("AD15").Value = 251
$ws.Range("AE15").Value = 12
$ws.Range("AF15").Value = 21
$ws.Range("AG15").Value = 13
$ws.Range("AH15").Value = 41
$ws.Range("AI15").Value = 34
$ws.Range("AJ15").Value = 41

# Row 16
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 2.9
$ws.Range("I16").Value = 2.3
$ws.Range("J16").Value = 1.08
$ws.Range("K16").Value = 8
$ws.Range("L16").Value = 1.36
$ws.Range("M16").Value = 3
$ws.Range("N16").Value = 2.15
$ws.Range("O16").Value = 1.67
$ws.Range("P16").Value = 1.44
$ws.Range("Q16").Value = 2.63
$ws.Range("R16").Value = 1.83
$ws.Range("S16").Value = 1.83
$ws.Range("T16").Value = 9
$ws.Range("U16").Value = 15
$ws.Range("V16").Value = 12
$ws.Range("W16").Value = 34
$ws.Range("X16").Value = 26
$ws.Range("Y16").Value = 41
$ws.Range("Z16").Value = 8
$ws.Range("AA16").Value = 6
$ws.Range("AB16").Value = 15
$ws.Range("AC16").Value = 51
$ws.Range("AD16").Value = 301
$ws.Range("AE16").Value = 7.5
$ws.Range("AF16").Value = 11
$ws.Range("AG16").Value = 10
$ws.Range("AH16").Value = 23
$ws.Range("AI16").Value = 21
$ws.Range("AJ16").Value = 34

# Row 17
$ws.Range("G17").Value = 1.8
$ws.Range("H17").Value = 3.4
$ws.Range("I17").Value = 3.8
$ws.Range("J17").Value = 1.06
$ws.Range("K17").Value = 9.5
$ws.Range("L17").Value = 1.33
$ws.Range("M17").Value = 3.25
$ws.Range("N17").Value = 2.03
$ws.Range("O17").Value = 1.78
$ws.Range("P17").Value = 1.4
$ws.Range("Q17").Value = 2.75
$ws.Range("R17").Value = 1.83
$ws.Range("S17").Value = 1.83
$ws.Range("T17").Value = 7
$ws.Range("U17").Value = 8.5
$ws.Range("V17").Value = 9
$ws.Range("W17").Value = 15
$ws.Range("X17").Value = 15
$ws.Range("Y17").Value = 29
$ws.Range("Z17").Value = 9.5
$ws.Range("AA17").Value = 7
$ws.Range("AB17").Value = 17
$ws.Range("AC17").Value = 51
$ws.Range("AD17").Value = 301
$ws.Range("AE17").Value = 11
$ws.Range("AF17").Value = 21
$ws.Range("AG17").Value = 13
$ws.Range("AH17").Value = 41
$ws.Range("AI17").Value = 34
$ws.Range("AJ17").Value = 41

# Row 18
$ws.Range("G18").Value = 1.95
$ws.Range("H18").Value = 3.3
$ws.Range("I18").Value = 3.3
$ws.Range("J18").Value = 1.05
$ws.Range("K18").Value = 11
$ws.Range("L18").Value = 1.29
$ws.Range("M18").Value = 3.5
$ws.Range("N18").Value = 1.95
$ws.Range("O18").Value = 1.85
$ws.Range("P18").Value = 1.36
$ws.Range("Q18").Value = 3
$ws.Range("R18").Value = 1.73
$ws.Range("S18").Value = 2
$ws.Range("T18").Value = 8
$ws.Range("U18").Value = 10
$ws.Range("V18").Value = 9
$ws.Range("W18").Value = 19
$ws.Range("X18").Value = 17
$ws.Range("Y18").Value = 26
$ws.Range("Z18").Value = 11
$ws.Range("AA18").Value = 6.5
$ws.Range("AB18").Value = 15
$ws.Range("AC18").Value = 51
$ws.Range("AD18").Value = 201
$ws.Range("AE18").Value = 11
$ws.Range("AF18").Value = 19
$ws.Range("AG18").Value = 12
$ws.Range("AH18").Value = 41
$ws.Range("AI18").Value = 29
$ws.Range("AJ18").Value = 34

# Row 19
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 3.3
$ws.Range("I19").Value = 3.4
$ws.Range("J19").Value = 1.07
$ws.Range("K19").Value = 9
$ws.Range("L19").Value = 1.33
$ws.Range("M19").Value = 3.25
$ws.Range("N19").Value = 2.08
$ws.Range("O19").Value = 1.73
$ws.Range("P19").Value = 1.44
$ws.Range("Q19").Value = 2.63
$ws.Range("R19").Value = 1.83
$ws.Range("S19").Value = 1.83
$ws.Range("T19").Value = 7
$ws.Range("U19").Value = 9
$ws.Range("V19").Value = 9
$ws.Range("W19").Value = 17
$ws.Range("X19").Value = 17
$ws.Range("Y19").Value = 29
$ws.Range("Z19").Value = 9
$ws.Range("AA19").Value = 6.5
$ws.Range("AB19").Value = 17
$ws.Range("AC19").Value = 51
$ws.Range("AD19").Value = 301
$ws.Range("AE19").Value = 9.5
$ws.Range("AF19").Value = 17
$ws.Range("AG19").Value = 13
$ws.Range("AH19").Value = 41
$ws.Range("AI19").Value = 29
$ws.Range("AJ19").Value = 41

# Row 24
$ws.Range("N24").Value = 2.05
$ws.Range("O24").Value = 1.62
$ws.Range("P24").Value = 1.39
$ws.Range("Q24").Value = 2.42
$ws.Range("T24").Value = 5
$ws.Range("U24").Value = 6.2
$ws.Range("V24").Value = 7
$ws.Range("W24").Value = 10.75
$ws.Range("X24").Value = 12
$ws.Range("Y24").Value = 25
$ws.Range("Z24").Value = 7.9
$ws.Range("AA24").Value = 5.7
$ws.Range("AB24").Value = 14
$ws.Range("AC24").Value = 70
$ws.Range("AD24").Value = 500
$ws.Range("AE24").Value = 9.25
$ws.Range("AF24").Value = 21
$ws.Range("AG24").Value = 12.5
$ws.Range("AH24").Value = 65
$ws.Range("AI24").Value = 40
$ws.Range("AJ24").Value = 45

# Row 25
$ws.Range("N25").Value = 2.07
$ws.Range("O25").Value = 1.6
$ws.Range("P25").Value = 1.4
$ws.Range("Q25").Value = 2.4
$ws.Range("T25").Value = 5.1
$ws.Range("U25").Value = 6.5
$ws.Range("V25").Value = 7
$ws.Range("W25").Value = 11.25
$ws.Range("X25").Value = 12.5
$ws.Range("Y25").Value = 24
$ws.Range("Z25").Value = 7.7
$ws.Range("AA25").Value = 5.5
$ws.Range("AB25").Value = 14
$ws.Range("AC25").Value = 70
$ws.Range("AD25").Value = 600
$ws.Range("AE25").Value = 8.75
$ws.Range("AF25").Value = 19
$ws.Range("AG25").Value = 12.5
$ws.Range("AH25").Value = 60
$ws.Range("AI25").Value = 40
$ws.Range("AJ25").Value = 45

# Row 26
$ws.Range("N26").Value = 2.05
$ws.Range("O26").Value = 1.62
$ws.Range("P26").Value = 1.4
$ws.Range("Q26").Value = 2.42
$ws.Range("T26").Value = 6.4
$ws.Range("U26").Value = 10
$ws.Range("V26").Value = 8
$ws.Range("W26").Value = 22
$ws.Range("X26").Value = 18
$ws.Range("Y26").Value = 26
$ws.Range("Z26").Value = 7.8
$ws.Range("AA26").Value = 5.1
$ws.Range("AB26").Value = 11.5
$ws.Range("AC26").Value = 50
$ws.Range("AD26").Value = 350
$ws.Range("AE26").Value = 6.7
$ws.Range("AF26").Value = 11
$ws.Range("AG26").Value = 8.25
$ws.Range("AH26").Value = 25
$ws.Range("AI26").Value = 19
$ws.Range("AJ26").Value = 26

# Row 34
$ws.Range("H34").Value = 3.3
$ws.Range("I34").Value = 3.9
$ws.Range("L34").Value = 1.4
$ws.Range("M34").Value = 2.75
$ws.Range("Z34").Value = 7.5
$ws.Range("AG34").Value = 15
$ws.Range("AI34").Value = 41

# Row 35
$ws.Range("J35").Value = 1.07
$ws.Range("K35").Value = 9
$ws.Range("N35").Value = 2.25
$ws.Range("O35").Value = 1.62

